$d = $word.ActiveDocument

$replacements = @(
    @("2023-12-21 Thursday", "2023-12-22 Friday"),
    @("45÷6=7, 3", "75÷8=9, 3"),
    @("24÷7=3, 3", "37÷8=4, 5"),
    @("59÷8=7, 3", "71÷4=17, 3"),
    @("63÷8=7, 7", "60÷2=30, 0"),
    @("87÷6=14, 3", "85÷2=42, 1"),
    @("61÷3=20, 1", "29÷3=9, 2"),
    @("14÷9=1, 5", "53÷6=8, 5"),
    @("32÷2=16, 0", "58÷7=8, 2"),
    @("41÷7=5, 6", "69÷6=11, 3"),
    @("12÷7=1, 5", "38÷5=7, 3"),
    @("77÷2=38, 1", "73÷5=14, 3"),
    @("53÷2=26, 1", "68÷2=34, 0"),
    @("31÷5=6, 1", "33÷7=4, 5"),
    @("70÷2=35, 0", "55÷5=11, 0"),
    @("85÷4=21, 1", "99÷7=14, 1"),
    @("87÷5=17, 2", "95÷4=23, 3"),
    @("28÷5=5, 3", "92÷2=46, 0"),
    @("38÷9=4, 2", "51÷7=7, 2"),
    @("57÷7=8, 1", "50÷7=7, 1"),
    @("16÷8=2, 0", "59÷2=29, 1"),
    @("77÷4=19, 1", "48÷4=12, 0"),
    @("69÷9=7, 6", "73÷6=12, 1"),
    @("66÷3=22, 0", "95÷7=13, 4"),
    @("72÷8=9, 0", "64÷8=8, 0"),
    @("22÷2=11, 0", "33÷9=3, 6")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    [void]$range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Host "Done applying replacements"
